# Requisitos.xlsx maintenance: two functional requirements (RF09 "Apresentar
# Historico de Pedidos do Cliente" and RF10 "Gerenciar Usuarios") and two
# non-functional requirements (old RNF02 "Seguranca de acesso" and old RNF03
# "Autenticacao de usuarios") were dropped from the requirements table.
# Deleting the sheet rows shifts everything below up (and shrinks the table /
# dimension / shared strings automatically); the surviving RNF rows then get
# their "Codigo" column renumbered so the RNF sequence stays contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("13:14").Delete()   # old RNF02 (Seguranca de acesso), RNF03 (Autenticacao de usuarios)
$ws.Rows("10:11").Delete()   # old RF09 (Apresentar Historico de Pedidos do Cliente), RF10 (Gerenciar Usuarios)

# Renumber the remaining RNF rows (old RNF04..RNF08 -> RNF02..RNF06).
$ws.Range("A11").Value = "RNF02"
$ws.Range("A12").Value = "RNF03"
$ws.Range("A13").Value = "RNF04"
$ws.Range("A14").Value = "RNF05"
$ws.Range("A15").Value = "RNF06"

# Restore the on-screen view: scrolled so row 9 is at the top, with F14 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("F14").Select()
